{"js": "// For every table in the document, set vertical alignment to \"bottom\" on\n// each cell of the first two rows (the repeated \"pageby\" header rows:\n// the treatment-group name row and the \"n\" row). This mirrors adding\n// <w:vAlign w:val=\"bottom\"/> to the <w:tcPr> of those cells in the OOXML.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\n// Only the first two rows of each table are affected.\nconst ROW_COUNT_TO_UPDATE = 2;\n\nconst rowsToTouch = [];\nfor (const table of tables.items) {\n  const items = table.rows.items;\n  const limit = Math.min(ROW_COUNT_TO_UPDATE, items.length);\n  for (let i = 0; i < limit; i++) {\n    const row = items[i];\n    row.cells.load(\"items\");\n    rowsToTouch.push(row);\n  }\n}\nawait context.sync();\n\nfor (const row of rowsToTouch) {\n  for (const cell of row.cells.items) {\n    cell.verticalAlignment = Word.VerticalAlignment.bottom;\n  }\n}\nawait context.sync();\n", "ps1": "# For every table in the document, set vertical alignment to \"bottom\" on\n# each cell of the first two rows (the repeated \"pageby\" header rows:\n# the treatment-group name row and the \"n\" row). This mirrors adding\n# <w:vAlign w:val=\"bottom\"/> to the <w:tcPr> of those cells in the OOXML.\n# wdCellAlignVerticalBottom = 3\n\n$d = $word.ActiveDocument\n\nforeach ($table in $d.Tables) {\n    $rowCount = [Math]::Min(2, $table.Rows.Count)\n    for ($r = 1; $r -le $rowCount; $r++) {\n        $row = $table.Rows.Item($r)\n        foreach ($cell in $row.Cells) {\n            $cell.VerticalAlignment = 3\n        }\n    }\n}\n"}
